$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Sexo" column (G2:G23): "Masculino" -> "Hombre", "Femenino" -> "Mujer"
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

# First pass: "Femenino" -> "Mujer" (processed first so it becomes the
# lower shared-string index, matching the saved workbook).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2
    if ($v -eq "Femenino") {
        $cell.Value = "Mujer"
    }
}

# Second pass: "Masculino" -> "Hombre"
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2
    if ($v -eq "Masculino") {
        $cell.Value = "Hombre"
    }
}

# Update the last active selection, as captured in the saved file.
$ws.Range("Q16").Select()
